$d = $word.ActiveDocument

$replacements = @(
    @{old = "Advanced Django Web Development 2016"; new = "Advanced Django Web Development (Oct 2018)"},
    @{old = "Advanced Express"; new = "Advanced Express (Oct 2019)"},
    @{old = "Advanced Node Js"; new = "Advanced Node Js (Oct 2019)"},
    @{old = "Advanced Php Debugging Techniques"; new = "Advanced Php Debugging Techniques (Jul 2018)"},
    @{old = "Advanced Python"; new = "Advanced Python (Oct 2018)"},
    @{old = "Building Apis In Php Using The Slim Micro Framework"; new = "Building Apis In Php Using The Slim Micro Framework (Mar 2019)"},
    @{old = "Design The Web Adding Dynamic Qr Codes"; new = "Design The Web Adding Dynamic Qr Codes (Nov 2018)"},
    @{old = "Designing Restful Apis"; new = "Designing Restful Apis (Oct 2018)"},
    @{old = "Extending Laravel With First Party Packages"; new = "Extending Laravel With First Party Packages (Nov 2018)"},
    @{old = "Learning Django"; new = "Learning Django (Oct 2018)"},
    @{old = "Learning Symfony 3"; new = "Learning Symfony 3 (Nov 2018)"},
    @{old = "Node Js Essential Training 3"; new = "Node Js Essential Training 3 (Oct 2019)"},
    @{old = "Node Js Microservices"; new = "Node Js Microservices (Oct 2019)"},
    @{old = "Node Js Security"; new = "Node Js Security (Nov 2019)"},
    @{old = "Pandas Essential Training"; new = "Pandas Essential Training (Oct 2018)"},
    @{old = "Php Design Patterns"; new = "Php Design Patterns (Oct 2018)"},
    @{old = "Php Testing Legacy Applications"; new = "Php Testing Legacy Applications (Nov 2018)"},
    @{old = "Python Advanced Design Patterns"; new = "Python Advanced Design Patterns (Nov 2018)"},
    @{old = "Typescript Essential Training"; new = "Typescript Essential Training (Dec 2018)"},
    @{old = "Working Remotely 2015"; new = "Working Remotely (Oct 2018)"}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
